# Update attendance computation results on Sheet1.
# Columns: D=Total Attendance Count, E=Real, F=Duplicate, G=Invalid, H=Absent
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: marked Invalid and Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Real attendance recorded
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Real attendance recorded
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent
$ws.Range("H6").Value = 1

# Row 7: Absent
$ws.Range("H7").Value = 1

# Row 8: Absent
$ws.Range("H8").Value = 1

# Row 9: Absent
$ws.Range("H9").Value = 1

# Row 10: Real attendance recorded
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: Real attendance recorded
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: Absent
$ws.Range("H12").Value = 1

# Row 13: Absent
$ws.Range("H13").Value = 1

# Row 14: Absent
$ws.Range("H14").Value = 1

# Row 15: Real attendance recorded
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

# Row 16: Absent
$ws.Range("H16").Value = 1

# Row 17: Absent
$ws.Range("H17").Value = 1

# Row 18: Absent
$ws.Range("H18").Value = 1
